$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sampling rate value (C4: 7 -> 64)
$ws.Range("C4").Value = 64

# Update CIC compensator decimation/interpolation coefficient (C10: 10 -> 1)
$ws.Range("C10").Value = 1

# Update CIC filter decimation/interpolation coefficient (C11: 256 -> 2560)
$ws.Range("C11").Value = 2560

# C12 and C15 are formulas (=C9/C11 and =C12*C4) and will recalc automatically

# Update the active selection to C11
$ws.Range("C11").Select()
